# Weekly data refresh: append 12 new match rows (115-126) to the "Main"
# sheet of the CS pre-match statistics workbook, then restore the
# worksheet's scroll/selection state to the new bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# Row -> column/value pairs, taken verbatim from the new sheetData rows.
# Missing columns in a row mean the source cell was blank (no stack count
# recorded, etc.) and are simply left untouched.
$newRows = @(
    @{ Row = 115; Cells = @{ A=874; B=5;     C=41;    D=20527;         F=22580; G=14769; H=18586; I=18695; J=18330; K=16498; L=19044; M=17640;         N=2;   O=-4;  P=0;   Q=-6;  R=-4;  S=2;  T=0  } },
    @{ Row = 116; Cells = @{ A=875; B=32;    C=5;     D=14338; E=18252; F=20216; G=14232; H=18846; I=22025; J=10625; K=22730;          M=15000; N=8;   O=-4;  P=1;   Q=-2;  R=-1;  S=2;  T=6  } },
    @{ Row = 117; Cells = @{ A=876; B=41;    C=11111; D=18061; E=14627; F=18522; G=20482; H=17114; I=18426; J=18363; K=18112; L=18067; M=18049; N=6;   O=-5;  P=3;   Q=-6;  R=0;   S=5;  T=-3 } },
    @{ Row = 118; Cells = @{ A=877; B=2111;  C=311;   D=18559; E=20000; F=15756; G=18154; H=17705; I=18245; J=18202; K=15421; L=13939; M=20900; N=-6;  O=5;   P=-8;  Q=-1;  R=0;   S=3;  T=7  } },
    @{ Row = 119; Cells = @{ A=878; B=2111;  C=311;   D=18474; E=15839; F=19767; G=17749; H=16970; I=18596; J=18292; K=20165; L=14509;          N=-10; O=3;   P=-1;  Q=5;   R=5;   S=1;  T=5  } },
    @{ Row = 120; Cells = @{ A=880; B=11111; C=221;   D=19026; E=18306; F=17717; G=15909; H=15857; I=15488; J=20000; K=14943; L=18780; M=15610; N=1;   O=8;   P=-1;  Q=-4;  R=-4;  S=0;  T=0  } },
    @{ Row = 121; Cells = @{ A=881; B=2111;  C=2111;  D=17951;          F=16925; G=13138; H=17518; I=17713; J=17184; K=16809; L=18369; M=13942; N=1;   O=-5;  P=-2;  Q=-6;  R=-2;  S=4;  T=4  } },
    @{ Row = 122; Cells = @{ A=882; B=311;   C=311;   D=21231;          F=21144; G=17668; H=17594; I=18277; J=15842; K=14340; L=20590; M=16594; N=1;   O=5;   P=1;   Q=1;   R=-4;  S=-1; T=-3 } },
    @{ Row = 123; Cells = @{ A=883; B=221;   C=311;   D=17183; E=16464; F=16677; G=16044; H=16197; I=18783; J=13952; K=16526; L=15815; M=15521; N=-2;  O=-15; P=1;   Q=6;   R=5;   S=1;  T=4  } },
    @{ Row = 124; Cells = @{ A=884; B=221;   C=5;     D=17511;          F=16085;                   I=17528; J=15323; K=17615;                   N=1;   O=2;   P=-4;  Q=-5;  R=-1;  S=-2; T=7  } },
    @{ Row = 125; Cells = @{ A=885; B=221;   C=2111;  D=16076; E=17131;          G=18786; H=10577; I=17772; J=14465; K=15662; L=15373; M=14999; N=1;   O=-2;  P=-6;  Q=3;   R=2;   S=-4; T=1  } },
    @{ Row = 126; Cells = @{ A=886; B=2111;  C=311;   D=17719; E=17840; F=17113; G=17042; H=16759; I=17764; J=17109; K=16757;          M=12117; N=1;   O=2;   P=2;   Q=2;   R=2;   S=-1; T=-3 } }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    foreach ($col in $entry.Cells.Keys) {
        $ws.Range("$col$r").Value = $entry.Cells[$col]
    }
}

# The sheet keeps its header row frozen; scroll the frozen pane down to the
# new last block of rows and leave the same relative selection (the cell
# just past the last data row) as the authored workbook shows.
$ws.Application.Goto($ws.Range("A115"), $true)
$ws.Range("O132").Select() | Out-Null
